# Update column G ("K") values on Sheet1 of the save-data workbook.
# This reflects regenerating the save data to use K (streak/strike count)
# instead of the old "Strike#" value, after recomputing std/mean and
# writing the new s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$gValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 2
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 2
    26 = 1
    27 = 2
    28 = 1
    29 = 1
    30 = 1
    31 = 2
    32 = 1
    33 = 2
    34 = 1
    35 = 0
    36 = 2
    37 = 1
    38 = 2
    39 = 0
    40 = 2
    41 = 0
    42 = 2
    43 = 0
    44 = 0
    45 = 0
    46 = 1
    47 = 0
    48 = 2
    49 = 1
    50 = 1
    51 = 0
    52 = 3
    54 = 1
    55 = 1
    56 = 2
    57 = 0
    58 = 2
    59 = 0
    60 = 1
    61 = 1
    62 = 1
    63 = 3
    64 = 1
    65 = 2
    66 = 0
    67 = 0
    68 = 2
    70 = 3
    72 = 2
    73 = 1
    75 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
